# Fruta / hortaliza, semanal
# Insert a new daily record as row 372, shifting the existing rows 372-472
# down to 373-473 (new dimension becomes A1:T473).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 372.
$ws.Rows.Item(372).Insert()

# Populate the new row with the new record's data.
$ws.Range("A372").Value = 9
$ws.Range("B372").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C372").Value = "Metropolitana"
$ws.Range("D372").Value = 44785
$ws.Range("E372").Value = 13
$ws.Range("F372").Value = "Fruta"
$ws.Range("G372").Value = 100108
$ws.Range("H372").Value = "Tropicales y subtropicales"
$ws.Range("I372").Value = 100108002
$ws.Range("J372").Value = "Mango"
$ws.Range("K372").Value = "Sin especificar"
$ws.Range("L372").Value = "Primera"
$ws.Range("M372").Value = 590
$ws.Range("N372").Value = 9500
$ws.Range("O372").Value = 10000
$ws.Range("P372").Value = 9754
$ws.Range("Q372").Value = "$/bandeja 4 kilos"
$ws.Range("R372").Value = "México"
$ws.Range("S372").Value = 2438
$ws.Range("T372").Value = 4
